$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 74: resale numbers update for 2024-01-19 13:02:50
# Columns A-D are textual (date/time/weekday/week) and must stay as text
# rather than being auto-converted to Excel date/number serials, so we
# force a Text number format while assigning them, then clear the
# formatting again so no stray style is left behind on the cells.
$textCells = "A74","B74","C74","D74"
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("A74").Value = "2024-01-19"
$ws.Range("B74").Value = "13:02:50"
$ws.Range("C74").Value = "Friday"
$ws.Range("D74").Value = "02"

foreach ($addr in $textCells) {
    $ws.Range($addr).ClearFormats()
}

$ws.Range("E74").Value = 137889
$ws.Range("F74").Value = 140445
$ws.Range("G74").Value = 171206
$ws.Range("H74").Value = 148815
$ws.Range("I74").Value = -1
$ws.Range("J74").Value = 121650
$ws.Range("K74").Value = 223453
$ws.Range("L74").Value = 254845
$ws.Range("M74").Value = 185299
$ws.Range("N74").Value = 110376
$ws.Range("O74").Value = 41337
$ws.Range("P74").Value = 30912
$ws.Range("Q74").Value = 73547
$ws.Range("R74").Value = -1
$ws.Range("S74").Value = 42494
$ws.Range("T74").Value = -1
